# Apply updated violent crime figures reflecting data for 2022-07-29.
# Each entry updates the 2022 (and, where applicable, 2021) year-to-date
# counts on the relevant worksheet.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = 'Citywide Totals'; Cell = 'I2'; Value = 4016 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I3'; Value = 4162 }
    @{ Sheet = 'Citywide Totals'; Cell = 'H4'; Value = 1667 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I4'; Value = 974 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I6'; Value = 4645 }
    @{ Sheet = 'Citywide Totals'; Cell = 'H7'; Value = 25976 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I7'; Value = 14179 }
    @{ Sheet = 'Uptown'; Cell = 'I2'; Value = 41 }
    @{ Sheet = 'Uptown'; Cell = 'I7'; Value = 159 }
    @{ Sheet = 'West Ridge'; Cell = 'I3'; Value = 43 }
    @{ Sheet = 'West Ridge'; Cell = 'I7'; Value = 151 }
    @{ Sheet = 'Fuller Park'; Cell = 'I3'; Value = 17 }
    @{ Sheet = 'Fuller Park'; Cell = 'I7'; Value = 48 }
    @{ Sheet = 'Grand Crossing'; Cell = 'I3'; Value = 146 }
    @{ Sheet = 'Grand Crossing'; Cell = 'I6'; Value = 125 }
    @{ Sheet = 'Grand Crossing'; Cell = 'I7'; Value = 456 }
    @{ Sheet = 'Woodlawn'; Cell = 'I3'; Value = 95 }
    @{ Sheet = 'Woodlawn'; Cell = 'I6'; Value = 75 }
    @{ Sheet = 'Woodlawn'; Cell = 'I7'; Value = 269 }
    @{ Sheet = 'North Lawndale'; Cell = 'I4'; Value = 29 }
    @{ Sheet = 'North Lawndale'; Cell = 'I6'; Value = 181 }
    @{ Sheet = 'North Lawndale'; Cell = 'I7'; Value = 552 }
    @{ Sheet = 'South Deering'; Cell = 'I2'; Value = 45 }
    @{ Sheet = 'South Deering'; Cell = 'I3'; Value = 40 }
    @{ Sheet = 'South Deering'; Cell = 'I6'; Value = 33 }
    @{ Sheet = 'South Deering'; Cell = 'I7'; Value = 124 }
    @{ Sheet = 'New City'; Cell = 'I6'; Value = 97 }
    @{ Sheet = 'New City'; Cell = 'I7'; Value = 319 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I7'; Value = 448 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I8'; Value = 852 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I10'; Value = 94 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I11'; Value = 216 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I12'; Value = 28 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I16'; Value = 38 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I18'; Value = 98 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I19'; Value = 392 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I22'; Value = 40 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I24'; Value = 36 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I25'; Value = 68 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I26'; Value = 22 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I29'; Value = 912 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I30'; Value = 48 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I33'; Value = 651 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I36'; Value = 196 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I37'; Value = 456 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I42'; Value = 489 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I44'; Value = 102 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I48'; Value = 190 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I51'; Value = 139 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I56'; Value = 16 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I57'; Value = 55 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I58'; Value = 10 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I61'; Value = 17 }
    @{ Sheet = 'By Neighborhood'; Cell = 'H63'; Value = 209 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I63'; Value = 53 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I65'; Value = 319 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I67'; Value = 552 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I69'; Value = 31 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I77'; Value = 78 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I79'; Value = 387 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I83'; Value = 290 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I84'; Value = 124 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I85'; Value = 639 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I86'; Value = 88 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I87'; Value = 29 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I88'; Value = 128 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I89'; Value = 159 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I90'; Value = 173 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I93'; Value = 84 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I94'; Value = 132 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I95'; Value = 230 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I96'; Value = 151 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I99'; Value = 269 }
    @{ Sheet = 'By Neighborhood'; Cell = 'H101'; Value = 25976 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I101'; Value = 14179 }
    @{ Sheet = 'South Chicago'; Cell = 'I3'; Value = 112 }
    @{ Sheet = 'South Chicago'; Cell = 'I7'; Value = 290 }
    @{ Sheet = 'West Pullman'; Cell = 'I3'; Value = 88 }
    @{ Sheet = 'West Pullman'; Cell = 'I4'; Value = 9 }
    @{ Sheet = 'West Pullman'; Cell = 'I7'; Value = 230 }
    @{ Sheet = 'Garfield Park'; Cell = 'I2'; Value = 154 }
    @{ Sheet = 'Garfield Park'; Cell = 'I3'; Value = 241 }
    @{ Sheet = 'Garfield Park'; Cell = 'I4'; Value = 30 }
    @{ Sheet = 'Garfield Park'; Cell = 'I6'; Value = 200 }
    @{ Sheet = 'Garfield Park'; Cell = 'I7'; Value = 651 }
    @{ Sheet = 'Englewood'; Cell = 'I2'; Value = 267 }
    @{ Sheet = 'Englewood'; Cell = 'I4'; Value = 46 }
    @{ Sheet = 'Englewood'; Cell = 'I6'; Value = 249 }
    @{ Sheet = 'Englewood'; Cell = 'I7'; Value = 912 }
    @{ Sheet = 'Chatham'; Cell = 'I2'; Value = 144 }
    @{ Sheet = 'Chatham'; Cell = 'I3'; Value = 113 }
    @{ Sheet = 'Chatham'; Cell = 'I6'; Value = 109 }
    @{ Sheet = 'Chatham'; Cell = 'I7'; Value = 392 }
    @{ Sheet = 'Irving Park'; Cell = 'I6'; Value = 32 }
    @{ Sheet = 'Irving Park'; Cell = 'I7'; Value = 102 }
    @{ Sheet = 'Lake View'; Cell = 'I6'; Value = 111 }
    @{ Sheet = 'Lake View'; Cell = 'I7'; Value = 190 }
    @{ Sheet = 'South Shore'; Cell = 'I3'; Value = 259 }
    @{ Sheet = 'South Shore'; Cell = 'I7'; Value = 639 }
    @{ Sheet = 'Humboldt Park'; Cell = 'I3'; Value = 167 }
    @{ Sheet = 'Humboldt Park'; Cell = 'I4'; Value = 41 }
    @{ Sheet = 'Humboldt Park'; Cell = 'I6'; Value = 132 }
    @{ Sheet = 'Humboldt Park'; Cell = 'I7'; Value = 489 }
    @{ Sheet = 'Avondale'; Cell = 'I6'; Value = 42 }
    @{ Sheet = 'Avondale'; Cell = 'I7'; Value = 94 }
    @{ Sheet = 'Dunning'; Cell = 'I2'; Value = 12 }
    @{ Sheet = 'Dunning'; Cell = 'I7'; Value = 36 }
    @{ Sheet = 'Norwood Park'; Cell = 'I2'; Value = 12 }
    @{ Sheet = 'Norwood Park'; Cell = 'I7'; Value = 31 }
    @{ Sheet = 'Roseland'; Cell = 'I2'; Value = 114 }
    @{ Sheet = 'Roseland'; Cell = 'I4'; Value = 25 }
    @{ Sheet = 'Roseland'; Cell = 'I6'; Value = 115 }
    @{ Sheet = 'Roseland'; Cell = 'I7'; Value = 387 }
    @{ Sheet = 'Calumet Heights'; Cell = 'I6'; Value = 42 }
    @{ Sheet = 'Calumet Heights'; Cell = 'I7'; Value = 98 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'I6'; Value = 58 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'I7'; Value = 196 }
    @{ Sheet = 'West Lawn'; Cell = 'I6'; Value = 36 }
    @{ Sheet = 'West Lawn'; Cell = 'I7'; Value = 84 }
    @{ Sheet = 'West Loop'; Cell = 'I3'; Value = 24 }
    @{ Sheet = 'West Loop'; Cell = 'I6'; Value = 73 }
    @{ Sheet = 'West Loop'; Cell = 'I7'; Value = 132 }
    @{ Sheet = 'East Side'; Cell = 'I2'; Value = 23 }
    @{ Sheet = 'East Side'; Cell = 'I7'; Value = 68 }
    @{ Sheet = 'East Village'; Cell = 'I2'; Value = 5 }
    @{ Sheet = 'East Village'; Cell = 'I7'; Value = 22 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'I2'; Value = 97 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'I7'; Value = 216 }
    @{ Sheet = 'United Center'; Cell = 'I2'; Value = 33 }
    @{ Sheet = 'United Center'; Cell = 'I7'; Value = 128 }
    @{ Sheet = 'Austin'; Cell = 'I2'; Value = 266 }
    @{ Sheet = 'Austin'; Cell = 'I3'; Value = 237 }
    @{ Sheet = 'Austin'; Cell = 'I6'; Value = 276 }
    @{ Sheet = 'Austin'; Cell = 'I7'; Value = 852 }
    @{ Sheet = 'Streeterville'; Cell = 'I2'; Value = 16 }
    @{ Sheet = 'Streeterville'; Cell = 'I7'; Value = 88 }
    @{ Sheet = 'Washington Heights'; Cell = 'I6'; Value = 59 }
    @{ Sheet = 'Washington Heights'; Cell = 'I7'; Value = 173 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'I3'; Value = 42 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'I7'; Value = 139 }
    @{ Sheet = 'Mckinley Park'; Cell = 'I3'; Value = 10 }
    @{ Sheet = 'Mckinley Park'; Cell = 'I7'; Value = 55 }
    @{ Sheet = 'Clearing'; Cell = 'I4'; Value = 3 }
    @{ Sheet = 'Clearing'; Cell = 'I7'; Value = 40 }
    @{ Sheet = 'Riverdale'; Cell = 'I6'; Value = 20 }
    @{ Sheet = 'Riverdale'; Cell = 'I7'; Value = 78 }
    @{ Sheet = 'Magnificent Mile'; Cell = 'I2'; Value = 5 }
    @{ Sheet = 'Magnificent Mile'; Cell = 'I7'; Value = 16 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'I3'; Value = 139 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'I7'; Value = 448 }
    @{ Sheet = 'Beverly'; Cell = 'I2'; Value = 7 }
    @{ Sheet = 'Beverly'; Cell = 'I7'; Value = 28 }
    @{ Sheet = 'Ukrainian Village'; Cell = 'I2'; Value = 3 }
    @{ Sheet = 'Ukrainian Village'; Cell = 'I7'; Value = 29 }
    @{ Sheet = 'Mount Greenwood'; Cell = 'I2'; Value = 7 }
    @{ Sheet = 'Mount Greenwood'; Cell = 'I7'; Value = 17 }
    @{ Sheet = 'Bucktown'; Cell = 'I6'; Value = 24 }
    @{ Sheet = 'Bucktown'; Cell = 'I7'; Value = 38 }
    @{ Sheet = 'Millenium Park'; Cell = 'I6'; Value = 4 }
    @{ Sheet = 'Millenium Park'; Cell = 'I7'; Value = 10 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
